# Apply updated odds values (FlashScore 2024-10-01 refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.45
$ws.Range("I2").Value = 1.93
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 2.27
$ws.Range("P2").Value = 3.9
$ws.Range("S2").Value = 1.32
$ws.Range("T2").Value = 3.1
$ws.Range("W2").Value = 13
$ws.Range("Y2").Value = 11.75
$ws.Range("Z2").Value = 45
$ws.Range("AA2").Value = 27
$ws.Range("AJ2").Value = 8.5
$ws.Range("AK2").Value = 17.5
$ws.Range("AN2").Value = 5.5
$ws.Range("AO2").Value = 18
$ws.Range("AP2").Value = 22
$ws.Range("AQ2").Value = 80
$ws.Range("AR2").Value = 100
$ws.Range("AT2").Value = 3.1
$ws.Range("AW2").Value = 4.05
$ws.Range("AX2").Value = 9.75
$ws.Range("BA2").Value = 55

# Row 3
$ws.Range("G3").Value = 1.93
$ws.Range("H3").Value = 3.85
$ws.Range("K3").Value = 2.4
$ws.Range("P3").Value = 4.75
$ws.Range("V3").Value = 2.57
$ws.Range("W3").Value = 11.75
$ws.Range("X3").Value = 12.5
$ws.Range("Z3").Value = 19
$ws.Range("AB3").Value = 17.5
$ws.Range("AH3").Value = 16
$ws.Range("AI3").Value = 22
$ws.Range("AJ3").Value = 11.75
$ws.Range("AN3").Value = 4.3
$ws.Range("AW3").Value = 5.7
$ws.Range("AX3").Value = 17
$ws.Range("AY3").Value = 19
$ws.Range("BA3").Value = 80

# Row 4
$ws.Range("G4").Value = 2.87
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 2.42
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 2.95
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 6.6
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 2.9
$ws.Range("Q4").Value = 2.07
$ws.Range("R4").Value = 1.7
$ws.Range("S4").Value = 1.39
$ws.Range("T4").Value = 2.77
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.91
$ws.Range("W4").Value = 8.5
$ws.Range("Z4").Value = 35
$ws.Range("AA4").Value = 25
$ws.Range("AB4").Value = 35
$ws.Range("AC4").Value = 6.6
$ws.Range("AD4").Value = 6
$ws.Range("AF4").Value = 70
$ws.Range("AG4").Value = 600
$ws.Range("AH4").Value = 7.4
$ws.Range("AK4").Value = 25
$ws.Range("AL4").Value = 21
$ws.Range("AM4").Value = 32
$ws.Range("AN4").Value = 4.85
$ws.Range("AO4").Value = 15
$ws.Range("AP4").Value = 21
$ws.Range("AR4").Value = 100
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.77
$ws.Range("AU4").Value = 6.7
$ws.Range("AV4").Value = 55
$ws.Range("AW4").Value = 4.35
$ws.Range("AY4").Value = 20
$ws.Range("AZ4").Value = 50
$ws.Range("BA4").Value = 80
$ws.Range("BB4").Value = 250

# Row 5
$ws.Range("I5").Value = 2.8
$ws.Range("L5").Value = 3.5
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("Q5").Value = 2.35
$ws.Range("R5").Value = 1.57
$ws.Range("U5").Value = 1.91
$ws.Range("V5").Value = 1.8
$ws.Range("W5").Value = 7.5
$ws.Range("AB5").Value = 41
$ws.Range("AC5").Value = 7
$ws.Range("AH5").Value = 7.5
$ws.Range("AL5").Value = 26
$ws.Range("AM5").Value = 41
$ws.Range("AP5").Value = 29
$ws.Range("AS5").Value = 251
$ws.Range("AX5").Value = 17
$ws.Range("AY5").Value = 29
$ws.Range("BB5").Value = 251

# Row 7
$ws.Range("BD7").Value = 126

# Row 8
$ws.Range("G8").Value = 2.15
$ws.Range("H8").Value = 3.25
$ws.Range("I8").Value = 3.5
$ws.Range("AJ8").Value = 12
$ws.Range("AY8").Value = 23
$ws.Range("BA8").Value = 67

# Row 9
$ws.Range("H9").Value = 3.9
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 11
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 3.75
$ws.Range("Q9").Value = 1.88
$ws.Range("R9").Value = 1.98
$ws.Range("S9").Value = 1.36
$ws.Range("T9").Value = 3
$ws.Range("U9").Value = 1.8
$ws.Range("V9").Value = 1.91
$ws.Range("X9").Value = 8
$ws.Range("AC9").Value = 12
$ws.Range("AE9").Value = 17
$ws.Range("AG9").Value = 251
$ws.Range("AL9").Value = 41
$ws.Range("AR9").Value = 51
$ws.Range("AT9").Value = 3

# Row 10
$ws.Range("O10").Value = 1.25
$ws.Range("P10").Value = 4
$ws.Range("Q10").Value = 1.82
$ws.Range("R10").Value = 1.92

# Row 11
$ws.Range("Q11").Value = 1.87
$ws.Range("R11").Value = 1.87

# Row 12
$ws.Range("Q12").Value = 1.97
$ws.Range("R12").Value = 1.77

# Row 13
$ws.Range("Q13").Value = 1.84
$ws.Range("R13").Value = 1.89
$ws.Range("BD13").Value = 151

# Row 14
$ws.Range("K14").Value = 2.2
$ws.Range("O14").Value = 1.25
$ws.Range("P14").Value = 3.75
$ws.Range("Q14").Value = 1.85
$ws.Range("R14").Value = 2
$ws.Range("X14").Value = 11
$ws.Range("AC14").Value = 11
$ws.Range("AL14").Value = 26

# Row 15
$ws.Range("Q15").Value = 1.75
$ws.Range("R15").Value = 2.05

# Row 16
$ws.Range("Q16").Value = 2
$ws.Range("R16").Value = 1.85

# Row 18
$ws.Range("G18").Value = 3.7
$ws.Range("I18").Value = 1.95
$ws.Range("J18").Value = 4.33
$ws.Range("L18").Value = 2.63
$ws.Range("M18").Value = 1.05
$ws.Range("N18").Value = 11
$ws.Range("Q18").Value = 2.03
$ws.Range("R18").Value = 1.83
$ws.Range("AA18").Value = 34
$ws.Range("AI18").Value = 9
$ws.Range("AK18").Value = 17
$ws.Range("AR18").Value = 101

# Row 19
$ws.Range("G19").Value = 2.7
$ws.Range("I19").Value = 2.45
$ws.Range("J19").Value = 3.2
$ws.Range("L19").Value = 3
$ws.Range("M19").Value = 1.04
$ws.Range("N19").Value = 13
$ws.Range("S19").Value = 1.33
$ws.Range("T19").Value = 3.25
$ws.Range("W19").Value = 11
$ws.Range("X19").Value = 15
$ws.Range("Y19").Value = 10
$ws.Range("AB19").Value = 23
$ws.Range("AC19").Value = 13
$ws.Range("AH19").Value = 11
$ws.Range("AL19").Value = 19
$ws.Range("AT19").Value = 3.25
$ws.Range("AW19").Value = 4.75
$ws.Range("AX19").Value = 13
$ws.Range("AY19").Value = 21
$ws.Range("BB19").Value = 126

# Row 21
$ws.Range("Q21").Value = 1.9
$ws.Range("R21").Value = 1.95

# Row 22
$ws.Range("G22").Value = 4.1
$ws.Range("H22").Value = 3.75
$ws.Range("I22").Value = 1.83
$ws.Range("L22").Value = 2.4
$ws.Range("S22").Value = 1.36
$ws.Range("T22").Value = 3
$ws.Range("W22").Value = 12
$ws.Range("AB22").Value = 41
$ws.Range("AI22").Value = 8.5
$ws.Range("AO22").Value = 23
$ws.Range("AQ22").Value = 81
$ws.Range("AR22").Value = 101
$ws.Range("AT22").Value = 3
$ws.Range("AW22").Value = 3.75
$ws.Range("AX22").Value = 9.5

# Row 23
$ws.Range("G23").Value = 2.6
$ws.Range("H23").Value = 3.3
$ws.Range("I23").Value = 2.75
$ws.Range("J23").Value = 3.4
$ws.Range("K23").Value = 2
$ws.Range("M23").Value = 1.08
$ws.Range("N23").Value = 7.5
$ws.Range("O23").Value = 1.44
$ws.Range("P23").Value = 2.63
$ws.Range("Q23").Value = 2.35
$ws.Range("R23").Value = 1.57
$ws.Range("S23").Value = 1.5
$ws.Range("T23").Value = 2.5
$ws.Range("U23").Value = 2
$ws.Range("V23").Value = 1.75
$ws.Range("W23").Value = 7
$ws.Range("Z23").Value = 26
$ws.Range("AA23").Value = 23
$ws.Range("AB23").Value = 41
$ws.Range("AC23").Value = 7.5
$ws.Range("AF23").Value = 67
$ws.Range("AG23").Value = 451
$ws.Range("AH23").Value = 7
$ws.Range("AI23").Value = 12
$ws.Range("AJ23").Value = 11
$ws.Range("AL23").Value = 26
$ws.Range("AM23").Value = 41
$ws.Range("AP23").Value = 29
$ws.Range("AS23").Value = 251
$ws.Range("AT23").Value = 2.5
$ws.Range("AU23").Value = 8.5
$ws.Range("AV23").Value = 67
$ws.Range("AY23").Value = 29
$ws.Range("BB23").Value = 251

# Row 24
$ws.Range("G24").Value = 2.4
$ws.Range("I24").Value = 2.8
$ws.Range("X24").Value = 13
$ws.Range("AH24").Value = 10

# Row 26
$ws.Range("G26").Value = 3.8
$ws.Range("I26").Value = 1.75
$ws.Range("AH26").Value = 10
$ws.Range("AN26").Value = 6
$ws.Range("AZ26").Value = 29
$ws.Range("BC26").Value = 351
